$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple +1 increments of Taxonsorteringsordning (column B) ---
$ws.Range("B2").Value  = 79003
$ws.Range("B3").Value  = 79245
$ws.Range("B4").Value  = 79245
$ws.Range("B5").Value  = 78648
$ws.Range("B6").Value  = 79245
$ws.Range("B8").Value  = 79864
$ws.Range("B9").Value  = 78911
$ws.Range("B10").Value = 79245
$ws.Range("B13").Value = 79245
$ws.Range("B14").Value = 79245
$ws.Range("B17").Value = 79245
$ws.Range("B18").Value = 79245
$ws.Range("B19").Value = 79245

# --- Rows 11 and 12 swap their species/location/time data (with B also +1) ---
$ws.Range("A11").Value  = 131033340
$ws.Range("B11").Value  = 79245
$ws.Range("E11").Value  = 6425
$ws.Range("F11").Value  = "Garnlav"
$ws.Range("G11").Value  = "Alectoria sarmentosa"
$ws.Range("H11").Value  = "(Ach.) Ach."
$ws.Range("Q11").Value  = 395817
$ws.Range("R11").Value  = 6804597
$ws.Range("Z11").Value  = "13:34"
$ws.Range("AB11").Value = "13:34"

$ws.Range("A12").Value  = 131033361
$ws.Range("B12").Value  = 78911
$ws.Range("E12").Value  = 353
$ws.Range("F12").Value  = "Dvärgbägarlav"
$ws.Range("G12").Value  = "Cladonia parasitica"
$ws.Range("H12").Value  = "(Hoffm.) Hoffm."
$ws.Range("Q12").Value  = 395806
$ws.Range("R12").Value  = 6804660
$ws.Range("Z12").Value  = "13:27"
$ws.Range("AB12").Value = "13:27"
